$d = $word.ActiveDocument

# 1. Remove one of the two consecutive empty paragraphs that sit right
#    before the first "List<TEntity> entityList = await _context.Set<TEntity>()"
#    code example (the one immediately preceding that paragraph).
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text
    if ($text -eq [char]13) {
        $next = $d.Paragraphs.Item($i + 1)
        $nextText = $next.Range.Text
        if ($nextText -eq "List<TEntity> entityList = await _context.Set<TEntity>()" + [char]13) {
            $p.Range.Delete()
            break
        }
    }
}

# 2. Update the "Ignore Auto Includes" paragraph text.
$d.Content.Find.Execute(
    "generate olan zaman IgnoreAutoIncludes method t",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "generate olan zaman,IQueryable-a IgnoreAutoIncludes method t",
    2
)
